$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)

# Locate the table shape on the slide (named "Table 4" in the deck).
$tblShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tblShape = $shp
    }
}

$tbl = $tblShape.Table

# Rows 2-6 hold the "Selectie" column (column 2) cells whose runs should no
# longer be forced bold (explicit b="0" is added to every run/endParaRPr).
for ($r = 2; $r -le 6; $r++) {
    $cell = $tbl.Rows($r).Cells(2)
    $tr = $cell.Shape.TextFrame.TextRange
    $tr.Font.Bold = 0
}
